$wb = $excel.ActiveWorkbook

# 1. Rename the "Include from ..." sheets to their new generic names.
$wsMondo = $wb.Worksheets.Item("Include from Mondo Disease On")
$wsMondo.Name = "Include #0"

$wsIcd10 = $wb.Worksheets.Item("Include from icd10-CA")
$wsIcd10.Name = "Include #1"

# 2. Update the Metadata sheet.
$ws = $wb.Worksheets.Item("Metadata")

# Date value changed.
$ws.Cells.Item(8, 2).Value = "2024-10-02T15:04:17+00:00"

# Contact value changed.
$ws.Cells.Item(10, 2).Value = "Ferlab.bio (http://example.org/example-publisher)"

# Insert a new "Jurisdiction" property row right after "Contact" (pushes
# Description/Purpose/Copyright/Immutable down by one row).
$ws.Range("A11").EntireRow.Insert()

$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = "'"

# Match the formatting of the other property rows for the newly inserted row
# (re-applied after the value write so the forced-text quote prefix doesn't
# leave behind its own one-off style).
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
